$d = $word.ActiveDocument

function Find-ParagraphLike($doc, [string]$pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Priority 2: Color coordinate portals using border/ center" -
#    split into several runs, wrapping "Color" and "center" in
#    proofErr spell-check markers (text itself is unchanged).
# ---------------------------------------------------------------------------
$pColor = Find-ParagraphLike $d "*Priority 2: Color coordinate portals using border/ center*"
if ($pColor -ne $null) {
    $xmlColor = '<w:p ' + $wNs + '>' +
        '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t xml:space="preserve">Priority 2: </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Color</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> coordinate portals using border/ </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>center</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'
    $pColor.Range.InsertXML($xmlColor) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) "Priority 3: Any other sounds in master list" - drop the stray
#    _GoBack bookmark that sits at the end of the paragraph.
# ---------------------------------------------------------------------------
$pSounds = Find-ParagraphLike $d "*Priority 3: Any other sounds in master list*"
if ($pSounds -ne $null) {
    $xmlSounds = '<w:p ' + $wNs + '>' +
        '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t>Priority 3: Any other sounds in master list</w:t></w:r>' +
        '</w:p>'
    $pSounds.Range.InsertXML($xmlSounds) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) After "Priority 2: Build the real level", add a new paragraph:
#    "Priority 3: Start implementing conversations" - with the
#    _GoBack bookmark now wrapping "Start implementing conversations".
# ---------------------------------------------------------------------------
$pBuild = Find-ParagraphLike $d "*Priority 2: Build the real level*"
if ($pBuild -ne $null) {
    # find the paragraph's 1-based index so the freshly-inserted paragraph
    # right after it can be looked up again post-insert
    $idx = 1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*Priority 2: Build the real level*") {
            $idx = $i
            break
        }
    }
    $pBuild.Range.InsertParagraphAfter() | Out-Null
    $newP = $d.Paragraphs.Item($idx + 1)
    $xmlNew = '<w:p ' + $wNs + '>' +
        '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t xml:space="preserve">Priority 3: </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Start implementing conversations</w:t></w:r>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
    $newP.Range.InsertXML($xmlNew) | Out-Null
}
